$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the BMI-related columns
$ws.Range("H1").Value = "Height"
$ws.Range("I1").Value = "Weight"

# Row 2 - existing data row, add height/weight
$ws.Range("H2").Value = -170
$ws.Range("I2").Value = 45

# Row 3 - existing data row, add height/weight
$ws.Range("H3").Value = 160
$ws.Range("I3").Value = 0

# Row 4 - existing data row, add height/weight
$ws.Range("H4").Value = 168468
$ws.Range("I4").Value = 54684

# Row 5 - new row with height/weight only
$ws.Range("H5").Value = 180
$ws.Range("I5").Value = 84

# Row 6 - new row with height/weight only
$ws.Range("H6").Value = 150.9
$ws.Range("I6").Value = 56.4

# Update selection to match final workbook state
$ws.Range("K4").Select()
